# Auto-generated Excel COM-interop script applying the scheduled-runner diff
# to the Adamantoise_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I43").Value = 2137.4167
$ws.Range("K43").Value = 2137.4167
$ws.Range("M43").Value = -2068.4167
$ws.Range("H43").Value = 3009.9333
$ws.Range("I51").Value = 4388.375
$ws.Range("M51").Value = -3904.375
$ws.Range("H51").Value = 3762.1035
$ws.Range("K51").Value = 4388.375
$ws.Range("H55").Value = 218.64285
$ws.Range("N55").Value = -628
$ws.Range("J55").Value = 200
$ws.Range("L55").Value = 200
$ws.Range("M58").Value = -1821
$ws.Range("H58").Value = 657
$ws.Range("K58").Value = 1971
$ws.Range("I58").Value = 657
$ws.Range("H76").Value = 4936.5
$ws.Range("H79").Value = 4936.5
$ws.Range("N112").Value = -10337.9999
$ws.Range("J112").Value = 2707.3333
$ws.Range("K112").Value = 5098.9998
$ws.Range("M112").Value = -3990.9998
$ws.Range("I112").Value = 1699.6666
$ws.Range("L112").Value = 8121.999899999999
$ws.Range("H112").Value = 2505.8
$ws.Range("M132").Value = -2132
$ws.Range("I132").Value = 1554
$ws.Range("K132").Value = 4662
$ws.Range("H132").Value = 1726.05
$ws.Range("H138").Value = 2622.45
$ws.Range("N138").Value = -19280
$ws.Range("I138").Value = 779.1177
$ws.Range("L138").Value = 9000
$ws.Range("K138").Value = 2337.3531
$ws.Range("M138").Value = 2802.6469
$ws.Range("J138").Value = 3000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2546.4546
$ws.Range("L2").Value = 3502.1667
$ws.Range("K2").Value = 1399.6
$ws.Range("I2").Value = 1399.6
$ws.Range("J2").Value = 3502.1667
$ws.Range("M2").Value = -1286.6
$ws.Range("N2").Value = -3728.1667
$ws.Range("K45").Value = 4547.3076
$ws.Range("M45").Value = -4170.3076
$ws.Range("I45").Value = 4547.3076
$ws.Range("H45").Value = 4157.4873
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("H110").Value = 2383.182
$ws.Range("M110").Value = 379.1111000000001
$ws.Range("I110").Value = 1665.8889
$ws.Range("K110").Value = 1665.8889
$ws.Range("I116").Value = 1399.6
$ws.Range("M116").Value = 894.4000000000001
$ws.Range("K116").Value = 1399.6
$ws.Range("L116").Value = 3502.1667
$ws.Range("N116").Value = -8090.1667
$ws.Range("J116").Value = 3502.1667
$ws.Range("H116").Value = 2546.4546
$ws.Range("M132").Value = -7049.119999999999
$ws.Range("I132").Value = 3193.04
$ws.Range("K132").Value = 9579.119999999999
$ws.Range("H132").Value = 3193.04
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 1399.6
$ws.Range("K3").Value = 1399.6
$ws.Range("N3").Value = -3730.1667
$ws.Range("H3").Value = 2546.4546
$ws.Range("J3").Value = 3502.1667
$ws.Range("M3").Value = -1285.6
$ws.Range("L3").Value = 3502.1667
$ws.Range("L63").Value = 100264
$ws.Range("H63").Value = 100264
$ws.Range("J63").Value = 100264
$ws.Range("N63").Value = -101636
$ws.Range("J66").Value = 100264
$ws.Range("L66").Value = 300792
$ws.Range("H66").Value = 100264
$ws.Range("N66").Value = -307656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J68").Value = 69967.25
$ws.Range("H68").Value = 69967.25
$ws.Range("N68").Value = -71465.25
$ws.Range("L68").Value = 69967.25
$ws.Range("L71").Value = 209901.75
$ws.Range("J71").Value = 69967.25
$ws.Range("N71").Value = -217389.75
$ws.Range("H71").Value = 69967.25
$ws.Range("L81").Value = 39000
$ws.Range("J81").Value = 39000
$ws.Range("N81").Value = -40996
$ws.Range("H81").Value = 39000
$ws.Range("N84").Value = -126984
$ws.Range("H84").Value = 39000
$ws.Range("L84").Value = 117000
$ws.Range("J84").Value = 39000
$ws.Range("K86").Value = 19544.455
$ws.Range("I86").Value = 19544.455
$ws.Range("M86").Value = -18421.455
$ws.Range("H86").Value = 22227.412
$ws.Range("N86").Value = -29392.166
$ws.Range("J86").Value = 27146.166
$ws.Range("L86").Value = 27146.166
$ws.Range("N87").Value = -74328
$ws.Range("L87").Value = 71956
$ws.Range("H87").Value = 71956
$ws.Range("J87").Value = 71956
$ws.Range("K89").Value = 97722.27500000001
$ws.Range("H89").Value = 22227.412
$ws.Range("L89").Value = 135730.83
$ws.Range("M89").Value = -92106.27500000001
$ws.Range("I89").Value = 19544.455
$ws.Range("N89").Value = -146962.83
$ws.Range("J89").Value = 27146.166
$ws.Range("J90").Value = 71956
$ws.Range("H90").Value = 71956
$ws.Range("L90").Value = 215868
$ws.Range("N90").Value = -227724
$ws.Range("J92").Value = 72309.5
$ws.Range("H92").Value = 72309.5
$ws.Range("N92").Value = -77301.5
$ws.Range("L92").Value = 72309.5
$ws.Range("H97").Value = 96180
$ws.Range("N97").Value = -98162
$ws.Range("L97").Value = 96180
$ws.Range("J97").Value = 96180
$ws.Range("K99").Value = 2422.6
$ws.Range("I99").Value = 2422.6
$ws.Range("H99").Value = 2685.5
$ws.Range("M99").Value = -924.5999999999999
$ws.Range("I126").Value = 2422.6
$ws.Range("H126").Value = 2685.5
$ws.Range("K126").Value = 7267.799999999999
$ws.Range("M126").Value = -4797.799999999999
$ws.Range("J132").Value = 1149.5
$ws.Range("N132").Value = -8508.5
$ws.Range("M132").Value = -7231.000100000001
$ws.Range("L132").Value = 3448.5
$ws.Range("I132").Value = 3253.6667
$ws.Range("K132").Value = 9761.000100000001
$ws.Range("H132").Value = 2412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I3").Value = 9311
$ws.Range("K3").Value = 27933
$ws.Range("H3").Value = 9311
$ws.Range("M3").Value = -27821
$ws.Range("J39").Value = 4994.8335
$ws.Range("N39").Value = -15572.5005
$ws.Range("H39").Value = 4671.5557
$ws.Range("L39").Value = 14984.5005
$ws.Range("L63").Value = 16665
$ws.Range("H63").Value = 5555
$ws.Range("J63").Value = 5555
$ws.Range("N63").Value = -18163
$ws.Range("J66").Value = 5555
$ws.Range("L66").Value = 49995
$ws.Range("H66").Value = 5555
$ws.Range("N66").Value = -57483
$ws.Range("L75").Value = 3000
$ws.Range("N75").Value = -4996
$ws.Range("J75").Value = 1000
$ws.Range("H75").Value = 1000
$ws.Range("J78").Value = 1000
$ws.Range("N78").Value = -18984
$ws.Range("H78").Value = 1000
$ws.Range("L78").Value = 9000
$ws.Range("H107").Value = 793
$ws.Range("L107").Value = 2969.3334
$ws.Range("J107").Value = 989.7778
$ws.Range("N107").Value = -6809.3334
$ws.Range("H113").Value = 1855
$ws.Range("M113").Value = -581
$ws.Range("I113").Value = 917
$ws.Range("K113").Value = 2751
$ws.Range("L121").Value = 5273.4
$ws.Range("K121").Value = 1495.28568
$ws.Range("H121").Value = 1023.1667
$ws.Range("N121").Value = -7893.4
$ws.Range("M121").Value = -185.28568
$ws.Range("J121").Value = 1757.8
$ws.Range("I121").Value = 498.42856
$ws.Range("N127").Value = -15920
$ws.Range("H127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("J127").Value = 2000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 413.41666
$ws.Range("L2").Value = 594
$ws.Range("K2").Value = 377.3
$ws.Range("I2").Value = 377.3
$ws.Range("J2").Value = 594
$ws.Range("M2").Value = -264.3
$ws.Range("N2").Value = -820
$ws.Range("J95").Value = 64933
$ws.Range("N95").Value = -70425
$ws.Range("L95").Value = 64933
$ws.Range("H95").Value = 64933
$ws.Range("K102").Value = 1391.75
$ws.Range("I102").Value = 1391.75
$ws.Range("H102").Value = 1500.4828
$ws.Range("M102").Value = 230.25
$ws.Range("J126").Value = 3326
$ws.Range("H126").Value = 2532.1924
$ws.Range("N126").Value = -14918
$ws.Range("L126").Value = 9978

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -2233.5
$ws.Range("I7").Value = 2345.5
$ws.Range("K7").Value = 2345.5
$ws.Range("H7").Value = 2353.2273
$ws.Range("H55").Value = 1511.742
$ws.Range("I55").Value = 1503.0555
$ws.Range("M55").Value = -1330.0555
$ws.Range("K55").Value = 1503.0555
$ws.Range("J64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9450
$ws.Range("H64").Value = 9000
$ws.Range("L67").Value = 9000
$ws.Range("H67").Value = 9000
$ws.Range("N67").Value = -10560
$ws.Range("J67").Value = 9000
$ws.Range("M100").Value = -1028.5714
$ws.Range("I100").Value = 1569.5714
$ws.Range("K100").Value = 1569.5714
$ws.Range("H100").Value = 1569.5714
$ws.Range("I126").Value = 2345.5
$ws.Range("H126").Value = 2353.2273
$ws.Range("K126").Value = 7036.5
$ws.Range("M126").Value = -4566.5
$ws.Range("L136").Value = 33999.999
$ws.Range("N136").Value = -39099.999
$ws.Range("M136").Value = -28386.375
$ws.Range("H136").Value = 10590.637
$ws.Range("I136").Value = 10312.125
$ws.Range("K136").Value = 30936.375
$ws.Range("J136").Value = 11333.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 746.73334
$ws.Range("I107").Value = 743.125
$ws.Range("L107").Value = 2252.5713
$ws.Range("M107").Value = -309.375
$ws.Range("K107").Value = 2229.375
$ws.Range("J107").Value = 750.8570999999999
$ws.Range("N107").Value = -6092.5713
$ws.Range("M136").Value = -3750
$ws.Range("H136").Value = 2753.9614
$ws.Range("I136").Value = 2100
$ws.Range("K136").Value = 6300

